$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Test_one / middle test / Last Name Test / 1
$ws.Range("A2").Value = "Test_one"
$ws.Range("A3").Value = "Test_two"

$ws.Range("B2").Value = "middle test"
$ws.Range("B3").Value = "Middle test two"

$ws.Range("C2").Value = "Last Name Test"
$ws.Range("C3").Value = "Last Name Test 2"

$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2

# Resize columns B and C to match the new (longer) best-fit text widths
$ws.Columns("B").ColumnWidth = 12.43
$ws.Columns("C").ColumnWidth = 13.43

# Update the active selection
$ws.Range("D4").Select()
